{"js": "// Add the author's email address after the \"Email Name of GitHub account to\n// gain access.\" sentence in the Assignments document.\n//\n// The underlying OOXML diff re-splits several runs around the word\n// \"GitHub\" (Word inserted w:proofErr spell/grammar-check markers while the\n// author was retyping) and relocates the \"_GoBack\" bookmark to the point of\n// the newest edit. Those are incidental artifacts of Word's live editor and\n// carry no semantic content change. The one actual content change is the\n// new text \" [archieoi[at]gmail.com]\" inserted right after\n// \"Email Name of GitHub account to gain access. \" (and before \"Under\n// folder...\"), which is what this script reproduces.\n\nconst body = context.document.body;\n\n// Find the exact sentence that precedes the newly-added email text.\nconst results = body.search(\"Email Name of GitHub account to gain access. \", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target sentence 'Email Name of GitHub account to gain access. '\");\n}\n\n// Insert the email address immediately after the located sentence.\nconst target = results.items[0];\ntarget.insertText(\" [archieoi[at]gmail.com]\", \"After\");\n\nawait context.sync();\n", "ps1": "# Add the author's email address after the \"Email Name of GitHub account to\n# gain access.\" sentence in the Assignments document.\n#\n# The underlying OOXML diff also re-splits several runs around the word\n# \"GitHub\" (Word inserted w:proofErr spell/grammar-check markers while the\n# author was retyping) and relocates the \"_GoBack\" bookmark to the point of\n# the newest edit. Those are incidental artifacts of Word's live editor and\n# carry no semantic content change. The one actual content change is the\n# new text \" [archieoi[at]gmail.com]\" inserted right after\n# \"Email Name of GitHub account to gain access. \" (and before \"Under\n# folder...\"), which is what this script reproduces.\n\n$d = $word.ActiveDocument\n\n$oldText = \"Email Name of GitHub account to gain access. \"\n$newText = \"Email Name of GitHub account to gain access.  [archieoi[at]gmail.com]\"\n\n# Find & replace the sentence in place (keeps the surrounding run\n# formatting \u2014 rFonts/color/sz \u2014 exactly as Word does on a live retype).\n$rng = $d.Content\n$found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\nif (-not $found) {\n    throw \"Could not find target sentence 'Email Name of GitHub account to gain access. '\"\n}\n\n$d.Save()\n"}
